$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value of 46081 (2026-02-28)
# for every data row (rows 2 through 514). Update it to 46082 (2026-03-01).
$ws.Range("C2:C514").Value = 46082
